# Applies the Aug 4 2023 cryptos-list refresh (prices / 1h-volume deltas,
# plus the Polkadot/WrappedEther and Stellar/Cosmos row swaps).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.135.29'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.03%  '

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.832.82'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.35%  '

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.03%  '

# Row 5: BNB
$ws.Range("E5").Value = '  +0.66%  '

# Row 6: XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6627'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.88%  '

# Row 7: USDC
$ws.Range("E7").Value = '  +0.04%  '

# Row 8: Dogecoin
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07414'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.66%  '

# Row 9: Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2934'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.05%  '

# Row 10: Solana
$ws.Range("E10").Value = '  -2.37%  '

# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07739'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.26%  '

# Row 12: WrappedEther (was Polkadot)
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.822.13'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.87%  '

# Row 13: Polkadot (was WrappedEther)
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.985'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.97%  '

# Row 14: Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6690'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.69%  '

# Row 15: Litecoin
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.87'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -5.64%  '

# Row 16: Uniswap
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.089'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.74%  '

# Row 17: ShibaInu
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008386'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.90%  '

# Row 18: WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.123.80'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.05%  '

# Row 19: BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '227.15'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.63%  '

# Row 20: Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.46'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.45%  '

# Row 21: Dai
$ws.Range("E21").Value = '  +0.17%  '

# Row 22: Chainlink
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.165'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.38%  '

# Row 23: BinanceUSD
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.05%  '

# Row 24: Monero
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '159.89'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.56%  '

# Row 25: Cosmos (was Stellar)
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.621'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.93%  '

# Row 26: Stellar (was Cosmos)
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1402'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.08%  '

# Row 27: EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.96'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.68%  '

# Row 28: PancakeSwap
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.512'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.45%  '

# Row 29: Filecoin
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.111'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -3.56%  '

# Row 30: InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.041'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.47%  '

# Row 31: Toncoin
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.193'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.02%  '

# Row 32: Hedera
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05353'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.25%  '

# Row 33: ImmutableX
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7563'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.27%  '

# Row 34: LidoDAOToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.872'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.01%  '

# Row 35: ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.135'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.15%  '

# Row 36: HuobiToken
$ws.Range("E36").Value = '  -0.46%  '

# Row 37: Maker
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.279.62'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.76%  '

# Row 38: VeChain
$ws.Range("E38").Value = '  -1.62%  '

# Row 39: MXToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.731'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.28%  '

# Row 40: TrustWalletToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9290'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.14%  '

# Row 41: XinFinNetwork
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.08940'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +15.72%  '

# Row 42: FraxShare
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.973'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.93%  '

# Row 43: PaxDollar
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.10%  '

# Row 44: Quant
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.85'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.71%  '

# Row 45: RocketPoolETH
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.964.72'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.92%  '

# Row 46: Mantle
$ws.Range("E46").Value = '  -0.44%  '

# Row 47: RenderToken
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.770'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.22%  '

# Row 48: BabyDogeCoin
$ws.Range("E48").Value = '  -1.44%  '

# Row 49: Aave
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '63.26'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.80%  '

# Row 50: Cronos
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05917'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.47%  '

# Row 51: Aptos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.796'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.26%  '
